# Auto-generated Excel COM-interop edit script
# Scheduled runner: refresh computed market-price / profit columns (H:N)
# for specific Leve rows across all job sheets, from latest market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 70.63636
$ws.Range("I38").Value = 70.63636
$ws.Range("K38").Value = 211.90908
$ws.Range("M38").Value = 160.09092

$ws.Range("H43").Value = 7748.3687
$ws.Range("I43").Value = 4165.4287
$ws.Range("J43").Value = 9838.416999999999
$ws.Range("K43").Value = 4165.4287
$ws.Range("L43").Value = 9838.416999999999
$ws.Range("M43").Value = -4096.4287
$ws.Range("N43").Value = -9976.416999999999

$ws.Range("H138").Value = 3848191
$ws.Range("I138").Value = 1121.3572
$ws.Range("J138").Value = 10873275
$ws.Range("K138").Value = 3364.0716
$ws.Range("L138").Value = 32619825
$ws.Range("M138").Value = 1775.9284
$ws.Range("N138").Value = -32630105

$ws.Range("H141").Value = 1932.0952
$ws.Range("I141").Value = 1731.3334
$ws.Range("J141").Value = 3136.6667
$ws.Range("K141").Value = 5194.0002
$ws.Range("L141").Value = 9410.000100000001
$ws.Range("M141").Value = -14.0002000000004
$ws.Range("N141").Value = -19770.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 668.75
$ws.Range("I5").Value = 262.5
$ws.Range("J5").Value = 1075
$ws.Range("K5").Value = 262.5
$ws.Range("L5").Value = 1075
$ws.Range("M5").Value = -150.5
$ws.Range("N5").Value = -1299

$ws.Range("H34").Value = 21624.8
$ws.Range("I34").Value = 4040
$ws.Range("J34").Value = 26021
$ws.Range("K34").Value = 4040
$ws.Range("L34").Value = 26021
$ws.Range("M34").Value = -3769
$ws.Range("N34").Value = -26563

$ws.Range("H45").Value = 5514.2856
$ws.Range("I45").Value = 11000
$ws.Range("J45").Value = 3320
$ws.Range("K45").Value = 11000
$ws.Range("L45").Value = 3320
$ws.Range("M45").Value = -10623
$ws.Range("N45").Value = -4074

$ws.Range("H88").Value = 2950
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 2933.3333
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 2933.3333
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -3745.3333

$ws.Range("H91").Value = 2950
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 2933.3333
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 2933.3333
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -5741.3333

$ws.Range("H122").Value = 29253
$ws.Range("I122").Value = 38004
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 114012
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -111562
$ws.Range("N122").Value = -13900

$ws.Range("H123").Value = 33974.3
$ws.Range("J123").Value = 33974.3
$ws.Range("L123").Value = 33974.3
$ws.Range("N123").Value = -43774.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 668.75
$ws.Range("I4").Value = 262.5
$ws.Range("J4").Value = 1075
$ws.Range("K4").Value = 262.5
$ws.Range("L4").Value = 1075
$ws.Range("M4").Value = -147.5
$ws.Range("N4").Value = -1305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 132750000
$ws.Range("J4").Value = 10333333
$ws.Range("L4").Value = 10333333
$ws.Range("N4").Value = -10333557

$ws.Range("H6").Value = 91366000
$ws.Range("I6").Value = 111668890
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 111668890
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -111668777
$ws.Range("N6").Value = -3226

$ws.Range("H7").Value = 201.90475
$ws.Range("I7").Value = 116.92857
$ws.Range("J7").Value = 371.85715
$ws.Range("K7").Value = 116.92857
$ws.Range("L7").Value = 371.85715
$ws.Range("M7").Value = -3.928569999999993
$ws.Range("N7").Value = -597.85715

$ws.Range("H17").Value = 508
$ws.Range("I17").Value = 508
$ws.Range("K17").Value = 508
$ws.Range("M17").Value = -334

$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10348

$ws.Range("H31").Value = 7411707.5
$ws.Range("I31").Value = 7013.76
$ws.Range("J31").Value = 16667575
$ws.Range("K31").Value = 7013.76
$ws.Range("L31").Value = 16667575
$ws.Range("M31").Value = -6718.76
$ws.Range("N31").Value = -16668165

$ws.Range("H34").Value = 7411707.5
$ws.Range("I34").Value = 7013.76
$ws.Range("J34").Value = 16667575
$ws.Range("K34").Value = 7013.76
$ws.Range("L34").Value = 16667575
$ws.Range("M34").Value = -6811.76
$ws.Range("N34").Value = -16667979

$ws.Range("H41").Value = 28766.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 28766.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 28766.25
$ws.Range("N41").Value = -29622.25
$ws.Range("M41").ClearContents()

$ws.Range("H70").Value = 32323.334
$ws.Range("J70").Value = 32323.334
$ws.Range("L70").Value = 32323.334
$ws.Range("N70").Value = -32953.334

$ws.Range("H73").Value = 32323.334
$ws.Range("J73").Value = 32323.334
$ws.Range("L73").Value = 32323.334
$ws.Range("N73").Value = -34507.334

$ws.Range("H107").Value = 434.92
$ws.Range("I107").Value = 244.76923
$ws.Range("J107").Value = 640.9167
$ws.Range("K107").Value = 244.76923
$ws.Range("L107").Value = 640.9167
$ws.Range("M107").Value = 1675.23077
$ws.Range("N107").Value = -4480.9167

$ws.Range("H122").Value = 4157
$ws.Range("I122").Value = 4671
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 14013
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -11563
$ws.Range("N122").Value = -16600

$ws.Range("H132").Value = 3052.6155
$ws.Range("I132").Value = 2575.3333
$ws.Range("J132").Value = 4126.5
$ws.Range("K132").Value = 7725.999899999999
$ws.Range("L132").Value = 12379.5
$ws.Range("M132").Value = -5195.999899999999
$ws.Range("N132").Value = -17439.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2233.3333
$ws.Range("J31").Value = 3000
$ws.Range("L31").Value = 9000
$ws.Range("N31").Value = -9576

$ws.Range("H38").Value = 167.88
$ws.Range("I38").Value = 114.666664
$ws.Range("J38").Value = 197.8125
$ws.Range("K38").Value = 343.999992
$ws.Range("L38").Value = 593.4375
$ws.Range("M38").Value = 3.00000799999998
$ws.Range("N38").Value = -1287.4375

$ws.Range("H75").Value = 4656.4736
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 5145.4707
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 15436.4121
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -17432.4121

$ws.Range("H78").Value = 4656.4736
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 5145.4707
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 46309.2363
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -56293.2363

$ws.Range("H129").Value = 2993.75
$ws.Range("I129").Value = 3298.8235
$ws.Range("J129").Value = 2648
$ws.Range("K129").Value = 9896.470499999999
$ws.Range("L129").Value = 7944
$ws.Range("M129").Value = -4896.470499999999
$ws.Range("N129").Value = -17944

$ws.Range("H131").Value = 741.215
$ws.Range("I131").Value = 423.88235
$ws.Range("J131").Value = 812.1974
$ws.Range("K131").Value = 1271.64705
$ws.Range("L131").Value = 2436.5922
$ws.Range("M131").Value = 3768.35295
$ws.Range("N131").Value = -12516.5922

$ws.Range("H134").Value = 3646.6667
$ws.Range("I134").Value = 2261.3333
$ws.Range("J134").Value = 5955.5557
$ws.Range("K134").Value = 6783.999899999999
$ws.Range("L134").Value = 17866.6671
$ws.Range("M134").Value = -1713.999899999999
$ws.Range("N134").Value = -28006.6671

$ws.Range("H137").Value = 10419413
$ws.Range("I137").Value = 11906115
$ws.Range("J137").Value = 12500
$ws.Range("K137").Value = 35718345
$ws.Range("L137").Value = 37500
$ws.Range("M137").Value = -35713245
$ws.Range("N137").Value = -47700

$ws.Range("H141").Value = 2047.1428
$ws.Range("I141").Value = 2047.1428
$ws.Range("K141").Value = 6141.428400000001
$ws.Range("M141").Value = -961.4284000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 97.15385000000001
$ws.Range("J2").Value = 190
$ws.Range("L2").Value = 190
$ws.Range("N2").Value = -416

$ws.Range("H15").Value = 19800
$ws.Range("J15").Value = 19800
$ws.Range("L15").Value = 19800
$ws.Range("N15").Value = -20376

$ws.Range("H81").Value = 19800
$ws.Range("J81").Value = 19800
$ws.Range("L81").Value = 19800
$ws.Range("N81").Value = -21796

$ws.Range("H84").Value = 19800
$ws.Range("J84").Value = 19800
$ws.Range("L84").Value = 59400
$ws.Range("N84").Value = -69384

$ws.Range("H122").Value = 33335332
$ws.Range("I122").Value = 66666664
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 199999992
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -199997542
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7655.1113
$ws.Range("I122").Value = 9875.5
$ws.Range("J122").Value = 5878.8
$ws.Range("K122").Value = 29626.5
$ws.Range("L122").Value = 17636.4
$ws.Range("M122").Value = -27176.5
$ws.Range("N122").Value = -22536.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 21436.111
$ws.Range("J86").Value = 21436.111
$ws.Range("L86").Value = 21436.111
$ws.Range("N86").Value = -23682.111

$ws.Range("H89").Value = 21436.111
$ws.Range("J89").Value = 21436.111
$ws.Range("L89").Value = 107180.555
$ws.Range("N89").Value = -118412.555

$ws.Range("H122").Value = 2452.7222
$ws.Range("I122").Value = 2289.9333
$ws.Range("J122").Value = 3266.6667
$ws.Range("K122").Value = 6869.7999
$ws.Range("L122").Value = 9800.000100000001
$ws.Range("M122").Value = -4419.7999
$ws.Range("N122").Value = -14700.0001

$ws.Range("H123").Value = 52300
$ws.Range("J123").Value = 52300
$ws.Range("L123").Value = 52300
$ws.Range("N123").Value = -62100

$ws.Range("H137").Value = 55700.8
$ws.Range("J137").Value = 55700.8
$ws.Range("L137").Value = 55700.8
$ws.Range("N137").Value = -65900.8
